$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy formatting (borders/fill/wrap) from an existing "standard" data row
# (style pattern s=7,7,8,7,7 across columns A:E) down onto the three new rows.
$ws.Range("A10:E10").Copy()
$ws.Range("A149:E151").PasteSpecial(-4122)

# Fill the new test cases column-by-column so the shared-string table grows
# in TCID, JIRA-ID, Description order (matches authoring order).
$ws.Range("A149").Value = "WAT166"
$ws.Range("A150").Value = "WAT167"
$ws.Range("A151").Value = "WAT168"

$ws.Range("B149").Value = "WAT-654"
$ws.Range("B150").Value = "WAT-655"
$ws.Range("B151").Value = "WAT-656"

$ws.Range("C149").Value = "Verify that user is able to sort author records/results using Sort by 'Relevance'"
$ws.Range("C150").Value = "Verify that user is able to sort author records/results using Sort by 'Publication years (newest first)'"
$ws.Range("C151").Value = "Verify that user is able to sort author records/results using Sort by 'Publication years (oldest first)'"

$ws.Range("D149").Value = "Y"
$ws.Range("D150").Value = "Y"
$ws.Range("D151").Value = "Y"

# Narrow the JIRA-ID column now that it no longer needs to fit the widest
# historical hyperlink text.
$ws.Columns.Item(2).ColumnWidth = 29.736979166666668

# Leave the sheet scrolled/selected near the newly-added rows.
$ws.Range("C155").Select() | Out-Null
